# Updates the cryptos list: refresh Price (column D) and Volume(1h) (column E)
# values for each coin row, matching the latest scraped figures.
# Equivalent to: "Updated cryptos list ... with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "97.935.28"
$ws.Range("E2").Value = "  +3.73%  "
$ws.Range("D3").Value = "3.352.11"
$ws.Range("E3").Value = "  +8.70%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'253.75"
$ws.Range("E5").Value = "  +6.97%  "
$ws.Range("D6").Value = "'622.94"
$ws.Range("E6").Value = "  +1.86%  "
$ws.Range("D7").Value = "'1.19"
$ws.Range("E7").Value = "  +7.19%  "
$ws.Range("D8").Value = "'0.384"
$ws.Range("E8").Value = "  +1.52%  "
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").Value = "3.347.97"
$ws.Range("E10").Value = "  +8.73%  "
$ws.Range("D11").Value = "'0.802"
$ws.Range("E11").Value = "  -1.30%  "
$ws.Range("E12").Value = "  +0.84%  "
$ws.Range("D13").Value = "97.790.17"
$ws.Range("E13").Value = "  +3.94%  "
$ws.Range("D14").Value = "'35.82"
$ws.Range("E14").Value = "  +5.29%  "
$ws.Range("E15").Value = "  +1.70%  "
$ws.Range("D16").Value = "3.967.85"
$ws.Range("E16").Value = "  +8.64%  "
$ws.Range("E17").Value = "  +2.37%  "
$ws.Range("D18").Value = "3.348.23"
$ws.Range("E18").Value = "  +8.42%  "
$ws.Range("D19").Value = "'3.61"
$ws.Range("E19").Value = "  +0.28%  "
$ws.Range("D20").Value = "'14.73"
$ws.Range("E20").Value = "  +1.89%  "
$ws.Range("D21").Value = "'478.73"
$ws.Range("E21").Value = "  +7.70%  "
$ws.Range("D22").Value = "'5.87"
$ws.Range("E22").Value = "  +2.19%  "
$ws.Range("D23").Value = "'0.0000207"
$ws.Range("E23").Value = "  +8.98%  "
$ws.Range("D24").Value = "'9.11"
$ws.Range("E24").Value = "  +3.04%  "
$ws.Range("D25").Value = "'5.68"
$ws.Range("E25").Value = "  +2.21%  "
$ws.Range("D26").Value = "'87.68"
$ws.Range("E26").Value = "  +3.45%  "
$ws.Range("D27").Value = "'11.93"
$ws.Range("E27").Value = "  -0.79%  "
$ws.Range("E28").Value = "  +10.25%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").Value = "'0.189"
$ws.Range("E30").Value = "  +5.66%  "
$ws.Range("D31").Value = "'0.252"
$ws.Range("E31").Value = "  -0.45%  "
$ws.Range("D32").Value = "'0.124"
$ws.Range("E32").Value = "  +0.69%  "
$ws.Range("E33").Value = "  +0.43%  "
$ws.Range("D34").Value = "'9.16"
$ws.Range("E34").Value = "  +1.49%  "
$ws.Range("D35").Value = "'27.16"
$ws.Range("E35").Value = "  +6.59%  "
$ws.Range("D36").Value = "'519.50"
$ws.Range("E36").Value = "  +7.08%  "
$ws.Range("D37").Value = "'0.151"
$ws.Range("E37").Value = "  -0.61%  "
$ws.Range("D38").Value = "'7.24"
$ws.Range("E38").Value = "  -6.13%  "
$ws.Range("D39").Value = "'1.93"
$ws.Range("E39").Value = "  +2.34%  "
$ws.Range("D40").Value = "'24.81"
$ws.Range("E40").Value = "  +3.04%  "
$ws.Range("D41").Value = "'0.448"
$ws.Range("E41").Value = "  +1.92%  "
$ws.Range("E42").Value = "  -0.60%  "
$ws.Range("D43").Value = "'3.67"
$ws.Range("E43").Value = "  -0.37%  "
$ws.Range("E44").Value = "  +16.83%  "
$ws.Range("D45").Value = "'3.22"
$ws.Range("E45").Value = "  +3.56%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").Value = "'161.00"
$ws.Range("E47").Value = "  -0.34%  "
$ws.Range("E48").Value = "  +5.25%  "
$ws.Range("D49").Value = "'45.51"
$ws.Range("E49").Value = "  +4.26%  "
$ws.Range("E50").Value = "  +6.09%  "
$ws.Range("D51").Value = "'4.51"
$ws.Range("E51").Value = "  +5.69%  "
